$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (shifts rows 3:22 down to 4:23)
$ws.Rows("3:3").Insert()

# Match the formatting used by the other date-label cells in column A
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the newly inserted row 3 with the missing quarter's data
$ws.Range("A3").Value = "2020-04-01 00:00:00_diff"
$ws.Range("B3").Value = 8.826710628892494
$ws.Range("C3").Value = -8.723943454208817
$ws.Range("D3").Value = -1.01454685455267
$ws.Range("E3").Value = 0.5683658063342414
$ws.Range("F3").Value = -2.397921186015015
$ws.Range("G3").Value = -0.1338900313505515
$ws.Range("H3").Value = -0.5850441862772902
